# Applies the "Add files via upload" revision to the Catering BI workbook.
#
# Scope of this script (the parts of the diff expressible through the Excel
# COM object model / this sandbox's iron_native shim):
#   1. Hide the daily-detail metric columns D:O on the "Base" sheet.
#   2. Fill in previously-blank KPI cells on row 87 (REC / RMSPII block) and
#      row 93 (RMSPII daily block) with the figures that were captured for
#      that date.
#   3. Update the sheet's frozen-pane scroll position / active selection to
#      reflect where the author was working when the file was saved.
#
# (Cosmetic, non-content package metadata such as fileVersion/build numbers,
# the xr:revisionPtr GUID, and the customXml SharePoint Flow bookkeeping
# parts are regenerated by Excel itself on every save and are not
# controllable - or meaningful - from the object model, so they are left
# alone here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Base")

# --- 1. Hide columns D:O -----------------------------------------------
$ws.Columns("D:O").Hidden = $true

# --- 2. Fill in the newly-reported KPI figures --------------------------

# Row 87 (REC, serial 45916 -> 2025-09-16): cycle-count / ILA / IRA /
# pending pallets / losses figures for that unit.
$ws.Range("G87").Value = 649
$ws.Range("H87").Value = 0.92710000000000004
$ws.Range("I87").Value = 0.98919999999999997
$ws.Range("J87").Value = 10
$ws.Range("K87").Value = 0

# Row 93 (RMSPII, serial 45917 -> 2025-09-17) is a hidden detail row;
# toggle it visible while writing so the engine doesn't bake an incidental
# row-height override into the saved file, then restore its hidden state.
$wasHidden93 = $ws.Rows(93).Hidden
$ws.Rows(93).Hidden = $false

$ws.Range("P93").Value = 207
$ws.Range("Q93").Value = 193
$ws.Range("R93").Value = 721064.78
$ws.Range("T93").Value = 1

# Z93 mirrors the "15/109" style ratio formulas already used elsewhere in
# column Z (e.g. Z77, Z85) - copy one of those cells' number format/font
# first so the new cell picks up the matching display style, then write the
# real formula over it.
$ws.Range("Z85").Copy() | Out-Null
$ws.Range("Z93").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("Z93").Formula = "=15/109"

$ws.Rows(93).Hidden = $wasHidden93

# --- 3. Update the saved view: frozen pane scroll + active selection ----
$ws.Range("H50").Select() | Out-Null
